$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. Values that are purely numeric-looking
# text (e.g. "4.40") are prefixed with a leading apostrophe so Excel keeps
# storing them as text (matching the original inline-string cell type) rather
# than silently converting them into a Double and losing formatting such as
# trailing zeros.
$changes = [ordered]@{
    'D2' = '43.195.97'
    'E2' = '  +0.93%  '
    'D3' = '2.353.98'
    'E3' = '  +2.24%  '
    'E4' = '  +0.03%  '
    'D5' = '''302.73'
    'E5' = '  +0.46%  '
    'D6' = '''95.67'
    'E6' = '  -0.56%  '
    'E7' = '  -1.42%  '
    'E8' = '  -0.04%  '
    'D9' = '''0.499'
    'E9' = '  +0.82%  '
    'D10' = '''34.19'
    'E10' = '  -1.76%  '
    'D11' = '''0.0788'
    'E11' = '  +0.01%  '
    'D12' = '''18.71'
    'E12' = '  -2.96%  '
    'E13' = '  +3.21%  '
    'E14' = '  -0.82%  '
    'D15' = '2.721.81'
    'E15' = '  +2.58%  '
    'D16' = '2.344.68'
    'E16' = '  +2.21%  '
    'E17' = '  +1.36%  '
    'D18' = '43.171.70'
    'E18' = '  +1.07%  '
    'D19' = '''12.22'
    'E19' = '  -1.00%  '
    'D20' = '''6.26'
    'E20' = '  +3.90%  '
    'E21' = '  -0.10%  '
    'D22' = '''68.16'
    'D23' = '''235.66'
    'E23' = '  +0.23%  '
    'D24' = '''2.23'
    'E24' = '  -1.41%  '
    'D26' = '''2.42'
    'E26' = '  +0.61%  '
    'D27' = '''24.55'
    'E27' = '  -0.35%  '
    'E28' = '  +14.63%  '
    'D29' = '''9.14'
    'E29' = '  +0.70%  '
    'D30' = '''31.44'
    'E30' = '  -2.52%  '
    'E31' = '  +0.04%  '
    'D32' = '''5.03'
    'E32' = '  +0.99%  '
    'D33' = '''0.0726'
    'E33' = '  +3.35%  '
    'D34' = '''17.19'
    'E34' = '  -1.64%  '
    'B35' = 'RenderToken'
    'C35' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D35' = '''4.40'
    'E35' = '  -0.78%  '
    'B36' = 'ARBITRUM'
    'C36' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D36' = '''1.84'
    'E36' = '  +4.68%  '
    'E37' = '  -0.75%  '
    'D38' = '''0.101'
    'E38' = '  +0.21%  '
    'B39' = 'LidoDAOToken'
    'C39' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D39' = '''2.76'
    'E39' = '  +1.52%  '
    'B40' = 'EnergySwap'
    'C40' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D40' = '''22.28'
    'E40' = '  +12.45%  '
    'E41' = '  -0.32%  '
    'D42' = '''110.67'
    'E42' = '  -32.91%  '
    'D43' = '1.944.10'
    'E43' = '  -1.44%  '
    'E44' = '  +0.25%  '
    'D45' = '''2.12'
    'E45' = '  +3.71%  '
    'D46' = '''9.41'
    'E46' = '  -10.06%  '
    'E47' = '  -0.86%  '
    'D48' = '2.585.01'
    'E48' = '  +2.32%  '
    'D49' = '''52.98'
    'E49' = '  -0.47%  '
    'E50' = '  -4.18%  '
    'D51' = '''72.21'
    'E51' = '  +0.82%  '
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
